{"js": "// Update benchmark stats table: replace the text of specific cells\n// (by row index, 0-based) in the single-column table with new values,\n// while preserving each cell's existing run formatting (font/size).\n// Rows 43-45 previously held multiple tab-separated numbers inside one\n// run; they are collapsed down to a single summary value each.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"97\",\n  4: \"0.00001\",\n  5: \"0.00044\",\n  6: \"0.00012\",\n  8: \"0.00017\",\n  9: \"0.00018\",\n  10: \"0.00027\",\n  11: \"0.01303\",\n  43: \"99.96\",\n  44: \"0.01\",\n  45: \"29\",\n};\n\nconst paragraphs = [];\nfor (const rowIndexStr of Object.keys(updates)) {\n  const rowIndex = parseInt(rowIndexStr, 10);\n  const cell = table.getCellOrNullObject(rowIndex, 0);\n  cell.load(\"body\");\n  paragraphs.push({ rowIndex, cell });\n}\nawait context.sync();\n\nfor (const { rowIndex, cell } of paragraphs) {\n  if (cell.isNullObject) {\n    continue;\n  }\n  cell.body.paragraphs.load(\"items\");\n}\nawait context.sync();\n\nfor (const { rowIndex, cell } of paragraphs) {\n  if (cell.isNullObject) {\n    continue;\n  }\n  const para = cell.body.paragraphs.items[0];\n  para.insertText(updates[rowIndex], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update benchmark stats table: replace the text of specific cells\n# (1-based row index) in the single-column table with new values,\n# while preserving each cell's existing run formatting (font/size).\n# Rows 44-46 previously held multiple tab-separated numbers inside one\n# run; they are collapsed down to a single summary value each.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$updates = [ordered]@{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"97\"\n    5  = \"0.00001\"\n    6  = \"0.00044\"\n    7  = \"0.00012\"\n    9  = \"0.00017\"\n    10 = \"0.00018\"\n    11 = \"0.00027\"\n    12 = \"0.01303\"\n    44 = \"99.96\"\n    45 = \"0.01\"\n    46 = \"29\"\n}\n\nforeach ($rowNum in $updates.Keys) {\n    $t.Cell($rowNum, 1).Range.Text = $updates[$rowNum]\n}\n"}
